# CIERRE 8 JUN 22
# Update the "VALES DE INSENTIVOS" sheet for May 2022 incentive payment,
# move the active tab from "ARQUITECTO" to "VALES DE INSENTIVOS", and
# re-centre (vertically) the signature line text in the stamp box.

$wb      = $excel.ActiveWorkbook
$wsVales = $wb.Worksheets.Item("VALES DE INSENTIVOS")

# Update the incentive-period label (ABRIL -> MAYO)
$wsVales.Range("A4").Value = "PAGO DE INCENTIVO DEL MES DE  MAYO   2022"

# Update the signature line, now padded with leading spaces, and centre it
# vertically within its merged box (C8:D9) along with the rest of the box.
$wsVales.Range("C8").Value = "                           PABLO BAEZ"
$wsVales.Range("C8").VerticalAlignment = -4108
$wsVales.Range("D8").VerticalAlignment = -4108
$wsVales.Range("C9").VerticalAlignment = -4108
$wsVales.Range("D9").VerticalAlignment = -4108

# Switch the active sheet/tab to "VALES DE INSENTIVOS" and move its
# selection to D14; "ARQUITECTO" loses its tabSelected flag as a result.
[void]$wsVales.Activate()
[void]$wsVales.Range("D14").Select()
